$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -746.7475945882428
$ws.Range("C2").Value = 15.68323435625896
$ws.Range("D2").Value = 1.30481513327601
$ws.Range("E2").Value = 34183

$ws.Range("B3").Value = -570.3977028459137
$ws.Range("C3").Value = 11.76935229067931
$ws.Range("D3").Value = 1.322748267898383
$ws.Range("E3").Value = 34182

$ws.Range("B4").Value = -427.7585592015432
$ws.Range("C4").Value = 10.55849741084228
$ws.Range("D4").Value = 1.323889246619446
$ws.Range("E4").Value = 34181
